$d = $word.ActiveDocument

$d.Content.Find.Execute(", and commitment purchasing strategies", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", and commitment purchasing strategies", 2)
